$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '25.897.62'
$ws.Range('E2').Value = '  -0.43%  '
$ws.Range('D3').Value = '1.641.39'
$ws.Range('E3').Value = '  +0.07%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.004'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.32%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '215.19'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.03%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.5056'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.05%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.004'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.41%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2570'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -0.15%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06395'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -0.59%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '19.62'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +0.84%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07801'
$ws.Range('D11').Style = 'Normal'
$ws.Range('B12').Value = 'Polkadot'
$ws.Range('C12').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '4.282'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +0.82%  '
$ws.Range('B13').Value = 'WrappedEther'
$ws.Range('C13').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D13').Value = '1.649.73'
$ws.Range('E13').Value = '  +0.50%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.5434'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -0.21%  '
$ws.Range('D15').Value = '0.0₅7877'
$ws.Range('E15').Value = '  -0.30%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '64.83'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +1.92%  '
$ws.Range('D17').Value = '25.937.73'
$ws.Range('E17').Value = '  -0.35%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '1.004'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -0.41%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '198.55'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -2.62%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '4.394'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +2.41%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '9.981'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.01%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.997'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.77%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '1.006'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.36%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.870'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -3.09%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '140.53'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -0.60%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.1144'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.82%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '6.851'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +1.83%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '15.75'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +0.08%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.244'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +0.25%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.04953'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -1.95%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '3.267'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +0.48%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.197'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +0.10%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.533'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -0.57%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '2.370'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +1.40%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.8937'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +0.38%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.607'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -0.94%  '
$ws.Range('D37').Value = '1.139.98'
$ws.Range('E37').Value = '  -0.47%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.5547'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -1.25%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01561'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -0.67%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.005'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -0.27%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '5.667'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +0.21%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.8215'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +1.50%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '99.37'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.51%  '
$ws.Range('D44').Value = '0.0₈121'
$ws.Range('E44').Value = '  +7.04%  '
$ws.Range('D45').Value = '1.777.85'
$ws.Range('E45').Value = '  +0.01%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.4521'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -0.13%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '55.38'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.82%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.004'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -0.33%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.05055'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +0.41%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.005'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -0.25%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.09518'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +2.14%  '
